$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.026462783288294
$ws.Cells.Item(2, 4).Value = 1.030756564931889
$ws.Cells.Item(2, 5).Value = 1.026685550906304
$ws.Cells.Item(2, 6).Value = 1.024979714447531
$ws.Cells.Item(2, 9).Value = 1.033926948399888
$ws.Cells.Item(2, 10).Value = 1.031626087283112
$ws.Cells.Item(2, 11).Value = 1.033566654395866
$ws.Cells.Item(2, 12).Value = 1.029507475098478
$ws.Cells.Item(2, 13).Value = 1.02780662802701
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.027365618354494
$ws.Cells.Item(3, 4).Value = 1.031420536549801
$ws.Cells.Item(3, 5).Value = 1.027450649110859
$ws.Cells.Item(3, 6).Value = 1.026516456930326
$ws.Cells.Item(3, 9).Value = 1.034157352540494
$ws.Cells.Item(3, 10).Value = 1.03216893292375
$ws.Cells.Item(3, 11).Value = 1.03403943531443
$ws.Cells.Item(3, 12).Value = 1.030080254816403
$ws.Cells.Item(3, 13).Value = 1.029148595215378
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.027949773309699
$ws.Cells.Item(4, 4).Value = 1.031850051738741
$ws.Cells.Item(4, 5).Value = 1.027946071304067
$ws.Cells.Item(4, 6).Value = 1.027510787652377
$ws.Cells.Item(4, 9).Value = 1.034305009878391
$ws.Cells.Item(4, 10).Value = 1.032519535329838
$ws.Cells.Item(4, 11).Value = 1.034344556731738
$ws.Cells.Item(4, 12).Value = 1.03045057565342
$ws.Cells.Item(4, 13).Value = 1.030016411982709
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028195342375544
$ws.Cells.Item(5, 4).Value = 1.032030590756867
$ws.Cells.Item(5, 5).Value = 1.028154430534757
$ws.Cells.Item(5, 6).Value = 1.027928798395322
$ws.Cells.Item(5, 9).Value = 1.034366742579107
$ws.Cells.Item(5, 10).Value = 1.032666771611646
$ws.Cells.Item(5, 11).Value = 1.034472638060112
$ws.Cells.Item(5, 12).Value = 1.030606184957996
$ws.Cells.Item(5, 13).Value = 1.030381120300875
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028236573929716
$ws.Cells.Item(6, 4).Value = 1.032060902321114
$ws.Cells.Item(6, 5).Value = 1.028189419858084
$ws.Cells.Item(6, 6).Value = 1.027998984030559
$ws.Cells.Item(6, 9).Value = 1.034377087684035
$ws.Cells.Item(6, 10).Value = 1.03269148402184
$ws.Cells.Item(6, 11).Value = 1.034494132226904
$ws.Cells.Item(6, 12).Value = 1.030632308120709
$ws.Cells.Item(6, 13).Value = 1.030442349403383
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.027953054653821
$ws.Cells.Item(7, 4).Value = 1.031852464225458
$ws.Cells.Item(7, 5).Value = 1.027948855082087
$ws.Cells.Item(7, 6).Value = 1.027516373149636
$ws.Cells.Item(7, 9).Value = 1.034305836099595
$ws.Cells.Item(7, 10).Value = 1.032521503325231
$ws.Cells.Item(7, 11).Value = 1.034346268914446
$ws.Cells.Item(7, 12).Value = 1.030452655202092
$ws.Cells.Item(7, 13).Value = 1.030021285706911
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.026767908744719
$ws.Cells.Item(8, 4).Value = 1.030980981428997
$ws.Cells.Item(8, 5).Value = 1.026944046053388
$ws.Cells.Item(8, 6).Value = 1.025499075819118
$ws.Cells.Item(8, 9).Value = 1.034005110431305
$ws.Cells.Item(8, 10).Value = 1.031809680043489
$ws.Cells.Item(8, 11).Value = 1.033726598258475
$ws.Cells.Item(8, 12).Value = 1.029701111834851
$ws.Cells.Item(8, 13).Value = 1.028260263395586
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.024679233058536
$ws.Cells.Item(9, 4).Value = 1.029444438179077
$ws.Cells.Item(9, 5).Value = 1.025176171849323
$ws.Cells.Item(9, 6).Value = 1.021943787557476
$ws.Cells.Item(9, 9).Value = 1.033464253587397
$ws.Cells.Item(9, 10).Value = 1.030550346084387
$ws.Cells.Item(9, 11).Value = 1.032628550395407
$ws.Cells.Item(9, 12).Value = 1.028374467392577
$ws.Cells.Item(9, 13).Value = 1.025152910786019
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.023286589072722
$ws.Cells.Item(10, 4).Value = 1.028419524776064
$ws.Cells.Item(10, 5).Value = 1.023999458685189
$ws.Cells.Item(10, 6).Value = 1.01957292554341
$ws.Cells.Item(10, 9).Value = 1.033096337683156
$ws.Cells.Item(10, 10).Value = 1.029707427672209
$ws.Cells.Item(10, 11).Value = 1.031892427943602
$ws.Cells.Item(10, 12).Value = 1.027488488944158
$ws.Cells.Item(10, 13).Value = 1.023078261464712
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.022683513287154
$ws.Cells.Item(11, 4).Value = 1.027975603018839
$ws.Cells.Item(11, 5).Value = 1.023490379221881
$ws.Cells.Item(11, 6).Value = 1.018546090449728
$ws.Cells.Item(11, 9).Value = 1.032935285785256
$ws.Cells.Item(11, 10).Value = 1.029341637763391
$ws.Cells.Item(11, 11).Value = 1.03157271150777
$ws.Cells.Item(11, 12).Value = 1.027104486094065
$ws.Cells.Item(11, 13).Value = 1.022179127920688
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022459496226217
$ws.Cells.Item(12, 4).Value = 1.027810692120232
$ws.Cells.Item(12, 5).Value = 1.02330135179908
$ws.Cells.Item(12, 6).Value = 1.018164636868392
$ws.Cells.Item(12, 9).Value = 1.032875202322423
$ws.Cells.Item(12, 10).Value = 1.029205646771389
$ws.Cells.Item(12, 11).Value = 1.031453808932986
$ws.Cells.Item(12, 12).Value = 1.026961795162861
$ws.Cells.Item(12, 13).Value = 1.021845025030465
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022507549003477
$ws.Cells.Item(13, 4).Value = 1.027846066900203
$ws.Cells.Item(13, 5).Value = 1.023341895767206
$ws.Cells.Item(13, 6).Value = 1.018246461919439
$ws.Cells.Item(13, 9).Value = 1.032888102266611
$ws.Cells.Item(13, 10).Value = 1.029234822746851
$ws.Cells.Item(13, 11).Value = 1.031479320527597
$ws.Cells.Item(13, 12).Value = 1.026992405346209
$ws.Cells.Item(13, 13).Value = 1.02191669695549
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.022664996127926
$ws.Cells.Item(14, 4).Value = 1.027961971808578
$ws.Cells.Item(14, 5).Value = 1.023474752779314
$ws.Cells.Item(14, 6).Value = 1.018514560254119
$ws.Cells.Item(14, 9).Value = 1.032930324603925
$ws.Cells.Item(14, 10).Value = 1.029330399160586
$ws.Cells.Item(14, 11).Value = 1.031562885954263
$ws.Cells.Item(14, 12).Value = 1.027092692341811
$ws.Cells.Item(14, 13).Value = 1.022151513443321
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.022762003453184
$ws.Cells.Item(15, 4).Value = 1.028033382200956
$ws.Cells.Item(15, 5).Value = 1.023556619309087
$ws.Cells.Item(15, 6).Value = 1.018679738838269
$ws.Cells.Item(15, 9).Value = 1.032956304515032
$ws.Cells.Item(15, 10).Value = 1.029389270996883
$ws.Cells.Item(15, 11).Value = 1.031614354075604
$ws.Cells.Item(15, 12).Value = 1.027154475163248
$ws.Cells.Item(15, 13).Value = 1.022296174964198
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.023326612128779
$ws.Cells.Item(16, 4).Value = 1.02844898376313
$ws.Cells.Item(16, 5).Value = 1.024033254031931
$ws.Cells.Item(16, 6).Value = 1.019641067733963
$ws.Cells.Item(16, 9).Value = 1.033106989474716
$ws.Cells.Item(16, 10).Value = 1.029731687071125
$ws.Cells.Item(16, 11).Value = 1.031913626035769
$ws.Cells.Item(16, 12).Value = 1.027513966202313
$ws.Cells.Item(16, 13).Value = 1.023137916788092
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.023680762319576
$ws.Cells.Item(17, 4).Value = 1.028709645909578
$ws.Cells.Item(17, 5).Value = 1.024332353912604
$ws.Cells.Item(17, 6).Value = 1.020244016614212
$ws.Cells.Item(17, 9).Value = 1.033201043875503
$ws.Cells.Item(17, 10).Value = 1.029946261238073
$ws.Cells.Item(17, 11).Value = 1.032101091699156
$ws.Cells.Item(17, 12).Value = 1.027739366937072
$ws.Cells.Item(17, 13).Value = 1.023665702400573
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.023887327147218
$ws.Cells.Item(18, 4).Value = 1.028861673391401
$ws.Cells.Item(18, 5).Value = 1.024506856726434
$ws.Cells.Item(18, 6).Value = 1.020595684097117
$ws.Cells.Item(18, 9).Value = 1.033255736144421
$ws.Cells.Item(18, 10).Value = 1.030071341469004
$ws.Cells.Item(18, 11).Value = 1.032210343593914
$ws.Cells.Item(18, 12).Value = 1.027870803822566
$ws.Cells.Item(18, 13).Value = 1.023973474283059
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.023957759572034
$ws.Cells.Item(19, 4).Value = 1.028913508723978
$ws.Cells.Item(19, 5).Value = 1.024566364940435
$ws.Cells.Item(19, 6).Value = 1.020715589998599
$ws.Cells.Item(19, 9).Value = 1.033274356284343
$ws.Cells.Item(19, 10).Value = 1.03011397749354
$ws.Cells.Item(19, 11).Value = 1.032247579776228
$ws.Cells.Item(19, 12).Value = 1.027915614368989
$ws.Cells.Item(19, 13).Value = 1.024078403710352
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.023642765872058
$ws.Cells.Item(20, 4).Value = 1.028681680600188
$ws.Cells.Item(20, 5).Value = 1.024300258883682
$ws.Cells.Item(20, 6).Value = 1.020179328263743
$ws.Cells.Item(20, 9).Value = 1.033190970112226
$ws.Cells.Item(20, 10).Value = 1.02992324745122
$ws.Cells.Item(20, 11).Value = 1.032080988093563
$ws.Cells.Item(20, 12).Value = 1.02771518724719
$ws.Cells.Item(20, 13).Value = 1.023609083951018
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.022618632091806
$ws.Cells.Item(21, 4).Value = 1.027927841210687
$ws.Cells.Item(21, 5).Value = 1.023435627837393
$ws.Cells.Item(21, 6).Value = 1.018435613163891
$ws.Cells.Item(21, 9).Value = 1.032917898399225
$ws.Cells.Item(21, 10).Value = 1.02930225761106
$ws.Cells.Item(21, 11).Value = 1.031538282038463
$ws.Cells.Item(21, 12).Value = 1.027063161843012
$ws.Cells.Item(21, 13).Value = 1.022082369321428
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.021974673000867
$ws.Cells.Item(22, 4).Value = 1.027453765168841
$ws.Cells.Item(22, 5).Value = 1.022892389920512
$ws.Cells.Item(22, 6).Value = 1.017339028380291
$ws.Cells.Item(22, 9).Value = 1.032744693777313
$ws.Cells.Item(22, 10).Value = 1.028911120416904
$ws.Cells.Item(22, 11).Value = 1.031196218560854
$ws.Cells.Item(22, 12).Value = 1.02665288844356
$ws.Cells.Item(22, 13).Value = 1.021121739534301
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.022316052245995
$ws.Cells.Item(23, 4).Value = 1.027705091833811
$ws.Cells.Item(23, 5).Value = 1.023180333497115
$ws.Cells.Item(23, 6).Value = 1.017920373601281
$ws.Cells.Item(23, 9).Value = 1.032836656269762
$ws.Cells.Item(23, 10).Value = 1.029118535654377
$ws.Cells.Item(23, 11).Value = 1.031377632703188
$ws.Cells.Item(23, 12).Value = 1.026870412334081
$ws.Cells.Item(23, 13).Value = 1.021631057821886
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.023659934846406
$ws.Cells.Item(24, 4).Value = 1.028694316957851
$ws.Cells.Item(24, 5).Value = 1.024314761113105
$ws.Cells.Item(24, 6).Value = 1.020208558212677
$ws.Cells.Item(24, 9).Value = 1.033195522531485
$ws.Cells.Item(24, 10).Value = 1.029933646629454
$ws.Cells.Item(24, 11).Value = 1.032090072336884
$ws.Cells.Item(24, 12).Value = 1.027726113118009
$ws.Cells.Item(24, 13).Value = 1.023634667626547
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.025219240338076
$ws.Cells.Item(25, 4).Value = 1.029841771489021
$ws.Cells.Item(25, 5).Value = 1.025632882754793
$ws.Cells.Item(25, 6).Value = 1.022863010937346
$ws.Cells.Item(25, 9).Value = 1.033605373362013
$ws.Cells.Item(25, 10).Value = 1.030876507187855
$ws.Cells.Item(25, 11).Value = 1.032913144282034
$ws.Cells.Item(25, 12).Value = 1.028717710997193
$ws.Cells.Item(25, 13).Value = 1.025956760369552
